$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"1"
$ws.Cells.Item(2, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2, 7).Value = [double]"0.06338433333333333"
$ws.Cells.Item(2, 8).Value = [double]"0.190153"
$ws.Cells.Item(2, 9).Value = [double]"0.001319770164420016"
$ws.Cells.Item(2, 10).Value = [double]"0.001319770164420016"
$ws.Cells.Item(2, 13).Value = [double]"788.1599833333333"
$ws.Cells.Item(2, 14).Value = [double]"2364.47995"
$ws.Cells.Item(2, 15).Value = [double]"0.8397951873720987"
$ws.Cells.Item(2, 16).Value = [double]"0.8397951873720988"
$ws.Cells.Item(2, 17).Value = [double]"49.95699510359444"
$ws.Cells.Item(2, 18).Value = [double]"449.61295593235"
$ws.Cells.Item(2, 19).Value = [double]"0.001108336632517213"
$ws.Cells.Item(2, 20).Value = [double]"0.001108336632517213"
$ws.Cells.Item(3, 5).Value = [double]"1"
$ws.Cells.Item(3, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(3, 7).Value = [double]"0.06338433333333333"
$ws.Cells.Item(3, 8).Value = [double]"0.190153"
$ws.Cells.Item(3, 9).Value = [double]"0.001319770164420016"
$ws.Cells.Item(3, 10).Value = [double]"0.001319770164420016"
$ws.Cells.Item(3, 15).Value = [double]"0.01890163353898316"
$ws.Cells.Item(3, 16).Value = [double]"0.01890163353898317"
$ws.Cells.Item(3, 17).Value = [double]"1.124403697896555"
$ws.Cells.Item(3, 18).Value = [double]"10.119633281069"
$ws.Cells.Item(3, 19).Value = [double]"2.49458120035507E-05"
$ws.Cells.Item(3, 20).Value = [double]"2.49458120035507E-05"
$ws.Cells.Item(4, 5).Value = [double]"1"
$ws.Cells.Item(4, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 7).Value = [double]"0.06338433333333333"
$ws.Cells.Item(4, 8).Value = [double]"0.190153"
$ws.Cells.Item(4, 9).Value = [double]"0.001319770164420016"
$ws.Cells.Item(4, 10).Value = [double]"0.001319770164420016"
$ws.Cells.Item(4, 13).Value = [double]"131.4690986666667"
$ws.Cells.Item(4, 14).Value = [double]"394.407296"
$ws.Cells.Item(4, 15).Value = [double]"0.1400821136357036"
$ws.Cells.Item(4, 16).Value = [double]"0.1400821136357036"
$ws.Cells.Item(4, 17).Value = [double]"8.333081172920888"
$ws.Cells.Item(4, 18).Value = [double]"74.99773055628799"
$ws.Cells.Item(4, 19).Value = [double]"0.0001848761941452959"
$ws.Cells.Item(4, 20).Value = [double]"0.0001848761941452959"
$ws.Cells.Item(5, 5).Value = [double]"1"
$ws.Cells.Item(5, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 7).Value = [double]"0.06338433333333333"
$ws.Cells.Item(5, 8).Value = [double]"0.190153"
$ws.Cells.Item(5, 9).Value = [double]"0.001319770164420016"
$ws.Cells.Item(5, 10).Value = [double]"0.001319770164420016"
$ws.Cells.Item(5, 13).Value = [double]"1.145987666666667"
$ws.Cells.Item(5, 14).Value = [double]"3.437963"
$ws.Cells.Item(5, 15).Value = [double]"0.001221065453214498"
$ws.Cells.Item(5, 16).Value = [double]"0.001221065453214498"
$ws.Cells.Item(5, 17).Value = [double]"0.07263766425988889"
$ws.Cells.Item(5, 18).Value = [double]"0.653738978339"
$ws.Cells.Item(5, 19).Value = [double]"1.611525753956498E-06"
$ws.Cells.Item(5, 20).Value = [double]"1.611525753956499E-06"
$ws.Cells.Item(6, 9).Value = [double]"0.4290986302854955"
$ws.Cells.Item(6, 10).Value = [double]"0.4290986302854955"
$ws.Cells.Item(6, 13).Value = [double]"788.1599833333333"
$ws.Cells.Item(6, 14).Value = [double]"2364.47995"
$ws.Cells.Item(6, 15).Value = [double]"0.8397951873720987"
$ws.Cells.Item(6, 16).Value = [double]"0.8397951873720988"
$ws.Cells.Item(6, 17).Value = [double]"16242.58431508945"
$ws.Cells.Item(6, 18).Value = [double]"146183.2588358051"
$ws.Cells.Item(6, 19).Value = [double]"0.3603549646217186"
$ws.Cells.Item(6, 20).Value = [double]"0.3603549646217186"
$ws.Cells.Item(7, 9).Value = [double]"0.4290986302854955"
$ws.Cells.Item(7, 10).Value = [double]"0.4290986302854955"
$ws.Cells.Item(7, 15).Value = [double]"0.01890163353898316"
$ws.Cells.Item(7, 16).Value = [double]"0.01890163353898317"
$ws.Cells.Item(7, 19).Value = [double]"0.008110665061736059"
$ws.Cells.Item(7, 20).Value = [double]"0.00811066506173606"
$ws.Cells.Item(8, 9).Value = [double]"0.4290986302854955"
$ws.Cells.Item(8, 10).Value = [double]"0.4290986302854955"
$ws.Cells.Item(8, 13).Value = [double]"131.4690986666667"
$ws.Cells.Item(8, 14).Value = [double]"394.407296"
$ws.Cells.Item(8, 15).Value = [double]"0.1400821136357036"
$ws.Cells.Item(8, 16).Value = [double]"0.1400821136357036"
$ws.Cells.Item(8, 17).Value = [double]"2709.345773799623"
$ws.Cells.Item(8, 18).Value = [double]"24384.11196419661"
$ws.Cells.Item(8, 19).Value = [double]"0.06010904308857755"
$ws.Cells.Item(8, 20).Value = [double]"0.06010904308857755"
$ws.Cells.Item(9, 9).Value = [double]"0.4290986302854955"
$ws.Cells.Item(9, 10).Value = [double]"0.4290986302854955"
$ws.Cells.Item(9, 13).Value = [double]"1.145987666666667"
$ws.Cells.Item(9, 14).Value = [double]"3.437963"
$ws.Cells.Item(9, 15).Value = [double]"0.001221065453214498"
$ws.Cells.Item(9, 16).Value = [double]"0.001221065453214498"
$ws.Cells.Item(9, 17).Value = [double]"23.61678046779711"
$ws.Cells.Item(9, 18).Value = [double]"212.551024210174"
$ws.Cells.Item(9, 19).Value = [double]"0.0005239575134632786"
$ws.Cells.Item(9, 20).Value = [double]"0.0005239575134632788"
$ws.Cells.Item(10, 5).Value = [double]"2"
$ws.Cells.Item(10, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10, 7).Value = [double]"0.1307036666666667"
$ws.Cells.Item(10, 8).Value = [double]"0.392111"
$ws.Cells.Item(10, 9).Value = [double]"0.002721473755033562"
$ws.Cells.Item(10, 10).Value = [double]"0.002721473755033562"
$ws.Cells.Item(10, 13).Value = [double]"788.1599833333333"
$ws.Cells.Item(10, 14).Value = [double]"2364.47995"
$ws.Cells.Item(10, 15).Value = [double]"0.8397951873720987"
$ws.Cells.Item(10, 16).Value = [double]"0.8397951873720988"
$ws.Cells.Item(10, 17).Value = [double]"103.0153997416055"
$ws.Cells.Item(10, 18).Value = [double]"927.1385976744499"
$ws.Cells.Item(10, 19).Value = [double]"0.002285480562036659"
$ws.Cells.Item(10, 20).Value = [double]"0.002285480562036659"
$ws.Cells.Item(11, 5).Value = [double]"2"
$ws.Cells.Item(11, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(11, 7).Value = [double]"0.1307036666666667"
$ws.Cells.Item(11, 8).Value = [double]"0.392111"
$ws.Cells.Item(11, 9).Value = [double]"0.002721473755033562"
$ws.Cells.Item(11, 10).Value = [double]"0.002721473755033562"
$ws.Cells.Item(11, 15).Value = [double]"0.01890163353898316"
$ws.Cells.Item(11, 16).Value = [double]"0.01890163353898317"
$ws.Cells.Item(11, 17).Value = [double]"2.318612161711445"
$ws.Cells.Item(11, 18).Value = [double]"20.867509455403"
$ws.Cells.Item(11, 19).Value = [double]"5.144029960360482E-05"
$ws.Cells.Item(11, 20).Value = [double]"5.144029960360483E-05"
$ws.Cells.Item(12, 5).Value = [double]"2"
$ws.Cells.Item(12, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(12, 7).Value = [double]"0.1307036666666667"
$ws.Cells.Item(12, 8).Value = [double]"0.392111"
$ws.Cells.Item(12, 9).Value = [double]"0.002721473755033562"
$ws.Cells.Item(12, 10).Value = [double]"0.002721473755033562"
$ws.Cells.Item(12, 13).Value = [double]"131.4690986666667"
$ws.Cells.Item(12, 14).Value = [double]"394.407296"
$ws.Cells.Item(12, 15).Value = [double]"0.1400821136357036"
$ws.Cells.Item(12, 16).Value = [double]"0.1400821136357036"
$ws.Cells.Item(12, 17).Value = [double]"17.18349324909511"
$ws.Cells.Item(12, 18).Value = [double]"154.651439241856"
$ws.Cells.Item(12, 19).Value = [double]"0.0003812297958091964"
$ws.Cells.Item(12, 20).Value = [double]"0.0003812297958091964"
$ws.Cells.Item(13, 5).Value = [double]"2"
$ws.Cells.Item(13, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(13, 7).Value = [double]"0.1307036666666667"
$ws.Cells.Item(13, 8).Value = [double]"0.392111"
$ws.Cells.Item(13, 9).Value = [double]"0.002721473755033562"
$ws.Cells.Item(13, 10).Value = [double]"0.002721473755033562"
$ws.Cells.Item(13, 13).Value = [double]"1.145987666666667"
$ws.Cells.Item(13, 14).Value = [double]"3.437963"
$ws.Cells.Item(13, 15).Value = [double]"0.001221065453214498"
$ws.Cells.Item(13, 16).Value = [double]"0.001221065453214498"
$ws.Cells.Item(13, 17).Value = [double]"0.1497847899881111"
$ws.Cells.Item(13, 18).Value = [double]"1.348063109893"
$ws.Cells.Item(13, 19).Value = [double]"3.323097584101416E-06"
$ws.Cells.Item(13, 20).Value = [double]"3.323097584101417E-06"
$ws.Cells.Item(14, 7).Value = [double]"27.224476"
$ws.Cells.Item(14, 8).Value = [double]"81.673428"
$ws.Cells.Item(14, 9).Value = [double]"0.566860125795051"
$ws.Cells.Item(14, 10).Value = [double]"0.566860125795051"
$ws.Cells.Item(14, 13).Value = [double]"788.1599833333333"
$ws.Cells.Item(14, 14).Value = [double]"2364.47995"
$ws.Cells.Item(14, 15).Value = [double]"0.8397951873720987"
$ws.Cells.Item(14, 16).Value = [double]"0.8397951873720988"
$ws.Cells.Item(14, 17).Value = [double]"21457.24255041873"
$ws.Cells.Item(14, 18).Value = [double]"193115.1829537686"
$ws.Cells.Item(14, 19).Value = [double]"0.4760464055558263"
$ws.Cells.Item(14, 20).Value = [double]"0.4760464055558263"
$ws.Cells.Item(15, 7).Value = [double]"27.224476"
$ws.Cells.Item(15, 8).Value = [double]"81.673428"
$ws.Cells.Item(15, 9).Value = [double]"0.566860125795051"
$ws.Cells.Item(15, 10).Value = [double]"0.566860125795051"
$ws.Cells.Item(15, 15).Value = [double]"0.01890163353898316"
$ws.Cells.Item(15, 16).Value = [double]"0.01890163353898317"
$ws.Cells.Item(15, 17).Value = [double]"482.9474394991827"
$ws.Cells.Item(15, 18).Value = [double]"4346.526955492644"
$ws.Cells.Item(15, 19).Value = [double]"0.01071458236563995"
$ws.Cells.Item(15, 20).Value = [double]"0.01071458236563995"
$ws.Cells.Item(16, 7).Value = [double]"27.224476"
$ws.Cells.Item(16, 8).Value = [double]"81.673428"
$ws.Cells.Item(16, 9).Value = [double]"0.566860125795051"
$ws.Cells.Item(16, 10).Value = [double]"0.566860125795051"
$ws.Cells.Item(16, 13).Value = [double]"131.4690986666667"
$ws.Cells.Item(16, 14).Value = [double]"394.407296"
$ws.Cells.Item(16, 15).Value = [double]"0.1400821136357036"
$ws.Cells.Item(16, 16).Value = [double]"0.1400821136357036"
$ws.Cells.Item(16, 17).Value = [double]"3579.177321392298"
$ws.Cells.Item(16, 18).Value = [double]"32212.59589253069"
$ws.Cells.Item(16, 19).Value = [double]"0.07940696455717157"
$ws.Cells.Item(16, 20).Value = [double]"0.07940696455717157"
$ws.Cells.Item(17, 7).Value = [double]"27.224476"
$ws.Cells.Item(17, 8).Value = [double]"81.673428"
$ws.Cells.Item(17, 9).Value = [double]"0.566860125795051"
$ws.Cells.Item(17, 10).Value = [double]"0.566860125795051"
$ws.Cells.Item(17, 13).Value = [double]"1.145987666666667"
$ws.Cells.Item(17, 14).Value = [double]"3.437963"
$ws.Cells.Item(17, 15).Value = [double]"0.001221065453214498"
$ws.Cells.Item(17, 16).Value = [double]"0.001221065453214498"
$ws.Cells.Item(17, 17).Value = [double]"31.19891372746266"
$ws.Cells.Item(17, 18).Value = [double]"280.790223547164"
$ws.Cells.Item(17, 19).Value = [double]"0.000692173316413161"
$ws.Cells.Item(17, 20).Value = [double]"0.0006921733164131611"
